$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark every checklist criterion as fulfilled (CALIFICACION = 1) except
# item 14 (row 26, which has no criterio text) which stays at 0.
$ws.Range("C13:C25").Value = 1
$ws.Range("C26").Value = 0
$ws.Range("C27").Value = 1

# Re-select the cell the author ended up on after filling the sheet in.
$ws.Range("C26").Select()
